$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fed")

# --- Step 1: set cell values (converts from shared-string text to numeric) ---
$ws.Range("Q3").Value = 5.15
$ws.Range("R3").Value = 5.15
$ws.Range("S3").Value = 5.15
$ws.Range("T3").Value = 5.85
$ws.Range("U3").Value = 6.55
$ws.Range("V3").Value = 7.25
$ws.Range("W3").Value = 7.25
$ws.Range("X3").Value = 7.25
$ws.Range("Y3").Value = 7.25
$ws.Range("Z3").Value = 7.25
$ws.Range("Q5").Value = 5.15
$ws.Range("R5").Value = 5.15
$ws.Range("Q6").Value = 5.15
$ws.Range("R6").Value = 5.15
$ws.Range("S6").Value = 6.25
$ws.Range("T6").Value = 6.25
$ws.Range("U6").Value = 6.25
$ws.Range("V6").Value = 6.25
$ws.Range("W6").Value = 6.25
$ws.Range("X6").Value = 6.25
$ws.Range("Y6").Value = 6.25
$ws.Range("Z6").Value = 6.25
$ws.Range("Q11").Value = 5.15
$ws.Range("Q12").Value = 5.15
$ws.Range("R12").Value = 5.15
$ws.Range("S12").Value = 5.15
$ws.Range("T12").Value = 5.15
$ws.Range("U12").Value = 5.15
$ws.Range("V12").Value = 5.15
$ws.Range("W12").Value = 5.15
$ws.Range("X12").Value = 5.15
$ws.Range("Y12").Value = 5.15
$ws.Range("Z12").Value = 5.15
$ws.Range("Q15").Value = 6.5
$ws.Range("R15").Value = 6.5
$ws.Range("S15").Value = 6.5
$ws.Range("T15").Value = 7.5
$ws.Range("U15").Value = 7.75
$ws.Range("V15").Value = 8
$ws.Range("W15").Value = 8.25
$ws.Range("X15").Value = 8.25
$ws.Range("Y15").Value = 8.25
$ws.Range("Z15").Value = 8.25
$ws.Range("Q16").Value = 5.15
$ws.Range("R16").Value = 5.15
$ws.Range("S16").Value = 5.15
$ws.Range("T16").Value = 5.85
$ws.Range("U16").Value = 6.55
$ws.Range("V16").Value = 7.25
$ws.Range("W16").Value = 7.25
$ws.Range("X16").Value = 7.25
$ws.Range("Y16").Value = 7.25
$ws.Range("Z16").Value = 7.25
$ws.Range("Q20").Value = 5.15
$ws.Range("R20").Value = 5.15
$ws.Range("S20").Value = 5.15
$ws.Range("T20").Value = 5.85
$ws.Range("U20").Value = 6.55
$ws.Range("V20").Value = 7.25
$ws.Range("W20").Value = 7.25
$ws.Range("X20").Value = 7.25
$ws.Range("Y20").Value = 7.25
$ws.Range("Z20").Value = 7.25
$ws.Range("Q24").Value = 5.15
$ws.Range("R24").Value = 5.15
$ws.Range("S24").Value = 6.95
$ws.Range("T24").Value = 7.15
$ws.Range("U24").Value = 7.4
$ws.Range("V24").Value = 7.4
$ws.Range("W24").Value = 7.4
$ws.Range("X24").Value = 7.4
$ws.Range("Y24").Value = 7.4
$ws.Range("Z24").Value = 8.15
$ws.Range("Q25").Value = 4.9
$ws.Range("R25").Value = 5.25
$ws.Range("S25").Value = 5.25
$ws.Range("T25").Value = 5.25
$ws.Range("U25").Value = 5.25
$ws.Range("V25").Value = 5.25
$ws.Range("W25").Value = 5.25
$ws.Range("X25").Value = 5.25
$ws.Range("Y25").Value = 5.25
$ws.Range("Z25").Value = 6.5
$ws.Range("Q26").Value = 5.15
$ws.Range("R26").Value = 5.15
$ws.Range("S26").Value = 5.15
$ws.Range("T26").Value = 5.85
$ws.Range("U26").Value = 6.55
$ws.Range("V26").Value = 7.25
$ws.Range("W26").Value = 7.25
$ws.Range("X26").Value = 7.25
$ws.Range("Y26").Value = 7.25
$ws.Range("Z26").Value = 7.25
$ws.Range("Q28").Value = 4
$ws.Range("R28").Value = 4
$ws.Range("S28").Value = 4
$ws.Range("T28").Value = 4
$ws.Range("U28").Value = 4
$ws.Range("V28").Value = 4
$ws.Range("W28").Value = 4
$ws.Range("X28").Value = 4
$ws.Range("Y28").Value = 4
$ws.Range("Z28").Value = 7.9
$ws.Range("Q29").Value = 5.15
$ws.Range("R29").Value = 5.15
$ws.Range("S29").Value = 5.15
$ws.Range("T29").Value = 5.85
$ws.Range("U29").Value = 6.55
$ws.Range("V29").Value = 7.25
$ws.Range("W29").Value = 7.25
$ws.Range("X29").Value = 7.25
$ws.Range("Y29").Value = 7.25
$ws.Range("Z29").Value = 7.25
$ws.Range("U30").Value = 6.55
$ws.Range("V30").Value = 6.55
$ws.Range("W30").Value = 7.25
$ws.Range("X30").Value = 7.25
$ws.Range("Y30").Value = 7.25
$ws.Range("Z30").Value = 7.25
$ws.Range("Q37").Value = 2.8
$ws.Range("R37").Value = 2.8
$ws.Range("Z37").Value = 7.25
$ws.Range("Q38").Value = 2
$ws.Range("R38").Value = 2
$ws.Range("S38").Value = 2
$ws.Range("T38").Value = 2
$ws.Range("U38").Value = 2
$ws.Range("V38").Value = 2
$ws.Range("W38").Value = 2
$ws.Range("X38").Value = 2
$ws.Range("Y38").Value = 2
$ws.Range("Z38").Value = 2
$ws.Range("Q42").Value = 5.15
$ws.Range("R42").Value = 5.15
$ws.Range("S42").Value = 5.15
$ws.Range("T42").Value = 5.85
$ws.Range("U42").Value = 6.55
$ws.Range("V42").Value = 7.25
$ws.Range("W42").Value = 7.25
$ws.Range("X42").Value = 7.25
$ws.Range("Y42").Value = 7.25
$ws.Range("Z42").Value = 7.25
$ws.Range("Q44").Value = 5.15
$ws.Range("R44").Value = 5.15
$ws.Range("S44").Value = 5.15
$ws.Range("T44").Value = 5.85
$ws.Range("U44").Value = 6.55
$ws.Range("V44").Value = 7.25
$ws.Range("W44").Value = 7.25
$ws.Range("X44").Value = 7.25
$ws.Range("Y44").Value = 7.25
$ws.Range("Z44").Value = 7.25
$ws.Range("Q47").Value = 7
$ws.Range("S47").Value = 7.53
$ws.Range("T47").Value = 7.68
$ws.Range("U47").Value = 8.06
$ws.Range("V47").Value = 8.06
$ws.Range("W47").Value = 8.15
$ws.Range("X47").Value = 8.46
$ws.Range("Y47").Value = 8.6
$ws.Range("Z47").Value = 8.73
$ws.Range("Q48").Value = 5.15
$ws.Range("R48").Value = 5.15
$ws.Range("S48").Value = 5.15
$ws.Range("T48").Value = 5.85
$ws.Range("U48").Value = 6.55
$ws.Range("V48").Value = 7.25
$ws.Range("W48").Value = 7.25
$ws.Range("X48").Value = 7.25
$ws.Range("Y48").Value = 7.25
$ws.Range("Z48").Value = 7.25
$ws.Range("Q50").Value = 5.15
$ws.Range("R50").Value = 5.15

# --- Step 2: apply style 2 (numFmt 40 "#,##0.00_);[Red](#,##0.00)" + black font) ---
$ws.Range("Q3").Font.Color = 0
$ws.Range("Q3").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R3").Font.Color = 0
$ws.Range("R3").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("S3").Font.Color = 0
$ws.Range("S3").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q5").Font.Color = 0
$ws.Range("Q5").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R5").Font.Color = 0
$ws.Range("R5").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q6").Font.Color = 0
$ws.Range("Q6").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R6").Font.Color = 0
$ws.Range("R6").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q11").Font.Color = 0
$ws.Range("Q11").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q12").Font.Color = 0
$ws.Range("Q12").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R12").Font.Color = 0
$ws.Range("R12").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q16").Font.Color = 0
$ws.Range("Q16").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R16").Font.Color = 0
$ws.Range("R16").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("S16").Font.Color = 0
$ws.Range("S16").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q20").Font.Color = 0
$ws.Range("Q20").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R20").Font.Color = 0
$ws.Range("R20").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("S20").Font.Color = 0
$ws.Range("S20").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q24").Font.Color = 0
$ws.Range("Q24").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R24").Font.Color = 0
$ws.Range("R24").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q26").Font.Color = 0
$ws.Range("Q26").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R26").Font.Color = 0
$ws.Range("R26").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("S26").Font.Color = 0
$ws.Range("S26").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q29").Font.Color = 0
$ws.Range("Q29").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R29").Font.Color = 0
$ws.Range("R29").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("S29").Font.Color = 0
$ws.Range("S29").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q42").Font.Color = 0
$ws.Range("Q42").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R42").Font.Color = 0
$ws.Range("R42").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("S42").Font.Color = 0
$ws.Range("S42").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q44").Font.Color = 0
$ws.Range("Q44").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R44").Font.Color = 0
$ws.Range("R44").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("S44").Font.Color = 0
$ws.Range("S44").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q48").Font.Color = 0
$ws.Range("Q48").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R48").Font.Color = 0
$ws.Range("R48").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("Q50").Font.Color = 0
$ws.Range("Q50").NumberFormat = "#,##0.00_);[Red](#,##0.00)"
$ws.Range("R50").Font.Color = 0
$ws.Range("R50").NumberFormat = "#,##0.00_);[Red](#,##0.00)"

# --- Step 3: apply style 3 (General numFmt + black font, new style) ---
$ws.Range("W3").Font.Color = 0
$ws.Range("X3").Font.Color = 0
$ws.Range("Y3").Font.Color = 0
$ws.Range("Z3").Font.Color = 0
$ws.Range("T16").Font.Color = 0
$ws.Range("U16").Font.Color = 0
$ws.Range("V16").Font.Color = 0
$ws.Range("T20").Font.Color = 0
$ws.Range("U20").Font.Color = 0
$ws.Range("V20").Font.Color = 0
$ws.Range("W20").Font.Color = 0
$ws.Range("X20").Font.Color = 0
$ws.Range("Y20").Font.Color = 0
$ws.Range("Z20").Font.Color = 0
$ws.Range("T26").Font.Color = 0
$ws.Range("U26").Font.Color = 0
$ws.Range("V26").Font.Color = 0
$ws.Range("W26").Font.Color = 0
$ws.Range("X26").Font.Color = 0
$ws.Range("Y26").Font.Color = 0
$ws.Range("Z26").Font.Color = 0
$ws.Range("T29").Font.Color = 0
$ws.Range("U29").Font.Color = 0
$ws.Range("W29").Font.Color = 0
$ws.Range("X29").Font.Color = 0
$ws.Range("Y29").Font.Color = 0
$ws.Range("Z29").Font.Color = 0
$ws.Range("U30").Font.Color = 0
$ws.Range("V30").Font.Color = 0
$ws.Range("W30").Font.Color = 0
$ws.Range("X30").Font.Color = 0
$ws.Range("Y30").Font.Color = 0
$ws.Range("Z30").Font.Color = 0
$ws.Range("T42").Font.Color = 0
$ws.Range("U42").Font.Color = 0
$ws.Range("V42").Font.Color = 0
$ws.Range("W42").Font.Color = 0
$ws.Range("X42").Font.Color = 0
$ws.Range("Y42").Font.Color = 0
$ws.Range("Z42").Font.Color = 0
$ws.Range("T44").Font.Color = 0
$ws.Range("U44").Font.Color = 0
$ws.Range("V44").Font.Color = 0
$ws.Range("W44").Font.Color = 0
$ws.Range("X44").Font.Color = 0
$ws.Range("Y44").Font.Color = 0
$ws.Range("Z44").Font.Color = 0
$ws.Range("T48").Font.Color = 0
$ws.Range("U48").Font.Color = 0
$ws.Range("V48").Font.Color = 0
$ws.Range("W48").Font.Color = 0
$ws.Range("X48").Font.Color = 0
$ws.Range("Y48").Font.Color = 0
$ws.Range("Z48").Font.Color = 0

# --- Step 4: update selection on the "fed" sheet ---
$ws.Activate()
$ws.Range("AB24").Select()
